$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-12 10:22:26", 0.0004),
    @("2023-12-12 10:23:00", 0.0022),
    @("2023-12-12 10:23:20", 0.0008),
    @("2023-12-12 10:23:35", 0.001),
    @("2023-12-12 10:23:41", 0.0004)
)

$startRow = 223
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
